$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "54.337.59"
Set-TextValue "E2" "  +0.60%  "
Set-TextValue "D3" "2.276.40"
Set-TextValue "E3" "  +1.43%  "
Set-TextValue "E4" "  +0.44%  "
Set-TextValue "D5" "498.61"
Set-TextValue "E5" "  +2.26%  "
Set-TextValue "D6" "128.32"
Set-TextValue "E6" "  +2.27%  "
Set-TextValue "E7" "  +0.22%  "
Set-TextValue "D8" "0.529"
Set-TextValue "E8" "  +1.10%  "
Set-TextValue "E9" "  +3.95%  "
Set-TextValue "E10" "  +1.94%  "
Set-TextValue "E11" "  +3.91%  "
Set-TextValue "E12" "  +1.79%  "
Set-TextValue "D13" "2.681.29"
Set-TextValue "E13" "  +1.98%  "
Set-TextValue "D14" "22.41"
Set-TextValue "E14" "  +5.44%  "
Set-TextValue "D15" "54.266.57"
Set-TextValue "E15" "  +1.13%  "
Set-TextValue "D16" "0.0000130"
Set-TextValue "E16" "  +1.21%  "
Set-TextValue "D17" "2.291.55"
Set-TextValue "E17" "  +2.79%  "
Set-TextValue "D18" "10.21"
Set-TextValue "E18" "  +5.70%  "
Set-TextValue "D19" "4.13"
Set-TextValue "E19" "  +3.17%  "
Set-TextValue "D20" "305.20"
Set-TextValue "E20" "  +3.02%  "
Set-TextValue "D21" "6.43"
Set-TextValue "E21" "  +3.58%  "
Set-TextValue "D22" "0.999"
Set-TextValue "E22" "  +0.19%  "
Set-TextValue "D23" "61.99"
Set-TextValue "E23" "  -2.74%  "
Set-TextValue "E24" "  +0.00%  "
Set-TextValue "D25" "2.381.72"
Set-TextValue "E25" "  +2.93%  "
Set-TextValue "E26" "  +2.53%  "
Set-TextValue "E27" "  +3.67%  "
Set-TextValue "D28" "173.59"
Set-TextValue "E28" "  +6.76%  "
Set-TextValue "E29" "  +2.26%  "
Set-TextValue "D30" "0.0₃0687"
Set-TextValue "E30" "  +2.52%  "
Set-TextValue "E31" "  +2.09%  "
Set-TextValue "E32" "  +2.69%  "
Set-TextValue "E33" "  +0.13%  "
Set-TextValue "E34" "  +2.44%  "
Set-TextValue "E35" "  -0.28%  "
Set-TextValue "D36" "0.922"
Set-TextValue "E36" "  +10.58%  "
Set-TextValue "E37" "  +1.70%  "
Set-TextValue "D38" "3.75"
Set-TextValue "E38" "  +4.83%  "
Set-TextValue "E39" "  +1.82%  "
Set-TextValue "E40" "  +1.12%  "
Set-TextValue "D41" "1.42"
Set-TextValue "E41" "  +2.18%  "
Set-TextValue "E42" "  +2.66%  "
Set-TextValue "E43" "  +2.80%  "
Set-TextValue "D44" "126.12"
Set-TextValue "E44" "  -1.29%  "
Set-TextValue "E46" "  +3.80%  "
Set-TextValue "E47" "  +1.93%  "
Set-TextValue "D48" "240.63"
Set-TextValue "E48" "  +1.40%  "
Set-TextValue "E49" "  +1.05%  "
Set-TextValue "E50" "  +2.16%  "
Set-TextValue "D51" "10.77"
Set-TextValue "E51" "  +1.08%  "
